$d = $word.ActiveDocument

function Replace-Text($old, $new) {
    $d.Content.Find.Execute($old, $true, $false, $false, $false, $false, `
                             $true, 1, $false, $new, 2) | Out-Null
}

Replace-Text "2025-07-24 Thursday" "2025-07-25 Friday"

Replace-Text "965×2=" "892×3="
Replace-Text "318×3=" "984×4="
Replace-Text "880×4=" "638×2="
Replace-Text "613×4=" "949×3="
Replace-Text "480×8=" "427×9="

Replace-Text "936×6=" "835×2="
Replace-Text "568×9=" "228×9="
Replace-Text "606×9=" "333×9="
Replace-Text "933×5=" "679×6="
Replace-Text "262×5=" "590×3="

Replace-Text "306×6=" "902×9="
Replace-Text "679×9=" "967×5="
Replace-Text "847×9=" "603×9="
Replace-Text "181×2=" "160×4="
Replace-Text "783×3=" "114×8="

Replace-Text "749×2=" "719×8="
Replace-Text "140×8=" "749×9="
Replace-Text "702×9=" "337×7="
Replace-Text "266×6=" "555×3="
Replace-Text "523×6=" "151×5="

Replace-Text "865×8=" "298×9="
Replace-Text "174×9=" "990×6="
Replace-Text "376×2=" "583×2="
Replace-Text "740×2=" "925×7="
Replace-Text "923×3=" "286×3="
